$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update birth_year (column Q) -1 and age_y (column S) +1 for rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("Q$r").Value2 = $ws.Range("Q$r").Value2 - 1
    $ws.Range("S$r").Value2 = $ws.Range("S$r").Value2 + 1
}
